$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 24, shifting existing rows 24:56 down to 25:57.
$ws.Rows.Item(24).Insert()

# Populate the newly inserted row 24 with the new weekly record.
$ws.Range("A24").Value = 3
$ws.Range("B24").Value = 'Femacal de La Calera'
$ws.Range("C24").Value = 'Coquimbo'
$ws.Range("D24").Value = 44757
$ws.Range("E24").Value = 5
$ws.Range("F24").Value = 100112035
$ws.Range("G24").Value = 'Bruselas (repollito)'
$ws.Range("H24").Value = 'Sin especificar'
$ws.Range("I24").Value = 'Primera'
$ws.Range("J24").Value = 40
$ws.Range("K24").Value = 15000
$ws.Range("L24").Value = 15000
$ws.Range("M24").Value = 15000
$ws.Range("N24").Value = '$/malla 15 kilos'
$ws.Range("O24").Value = 'Provincia de Quillota'
$ws.Range("P24").Value = 1000
$ws.Range("Q24").Value = 15
$ws.Range("R24").Value = 'Hortaliza'
